$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new repo row (row 4)
$ws.Range("A4").Value = "https://github.com/repoaleatorio/repoaleatorio"
$ws.Range("B4").Value = "maria"
$ws.Range("C4").Value = "email@example.com"

# Apply the same "Hiperlink" style used by the other link cells (A2,C2,A3,C3)
$ws.Range("A4").Style = "Hiperlink"
$ws.Range("C4").Style = "Hiperlink"

# Add hyperlinks for the new row
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/repoaleatorio/repoaleatorio") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:email@example.com") | Out-Null

$ws.Range("C4").Select() | Out-Null
